$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update existing row 4 values
$ws.Range("B4").Value = "Fulano Sousa 1"
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 80

# Add new row 8
$ws.Range("A8").Value = "AE123456789BR"
$ws.Range("B8").Value = "Fulano da Silva 5"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 13610720
$ws.Range("E8").Value = 3.9

# Add new row 9
$ws.Range("A9").Value = "AF123456789BR"
$ws.Range("B9").Value = "Fulano da Silva 6"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 19050280
$ws.Range("E9").Value = 5.6
